$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SEA")
$ws.Activate()

# --- Row 14 : fill in the Estimasi (N/O) dates + Note (P) for the "F1250" shipment ---
# Copy number-format/border styling from the row-8 pattern (same visual group: plain
# text row with date-formatted Estimasi columns + a note cell) onto row 14's N:O and P.
$ws.Range("N8:O8").Copy() | Out-Null
$ws.Range("N14:O14").PasteSpecial(-4122) | Out-Null
$ws.Range("P8").Copy() | Out-Null
$ws.Range("P14").PasteSpecial(-4122) | Out-Null

$ws.Range("N14").Value = 45672
$ws.Range("O14").Value = 45700
$ws.Range("P14").Value = "F1252"

# --- Row 16 : fill in the Estimasi (N/O) dates + Note (P) for the "F1270" shipment ---
# Row 16 heads a 16:22 merged block, so copy styling from row 5 (same merged-header
# shape) for the date columns, and from the already-merged G16 cell for the note column.
$ws.Range("N5:O5").Copy() | Out-Null
$ws.Range("N16:O16").PasteSpecial(-4122) | Out-Null
$ws.Range("G16").Copy() | Out-Null
$ws.Range("P16").PasteSpecial(-4122) | Out-Null

$ws.Range("N16").Value = 45680
$ws.Range("O16").Value = 45712
$ws.Range("P16").Value = "F1270"

# Merge the new Estimasi/Note columns across the row 16:22 block, matching the existing
# G16:G22 / I16:I22 / B16:B22 / C16:C22 / D16:D22 / E16:E22 merges already on that block.
$ws.Range("N16:N22").Merge() | Out-Null
$ws.Range("O16:O22").Merge() | Out-Null
$ws.Range("P16:P22").Merge() | Out-Null

# Restore the cursor/selection to the block the author was last working in.
$ws.Range("B16:B22").Select() | Out-Null
